$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(569).Delete()
